$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Materials" index numbers in column A for rows 2-4
$ws.Range("A2").Value = 9
$ws.Range("A3").Value = 10
$ws.Range("A4").Value = 11

# Remove the last two rows (formerly "D" and "E" materials), which
# shrinks the sheet dimension down to A1:L4 and drops the now-unused
# "D"/"E" shared strings.
$ws.Range("A5:L6").EntireRow.Delete()
